# Fix output.xlsx bimode global and choice columns.
#
# For every row whose "Branch Predictor" (column C) is "BiMode", the
# original data had the Local/Global/Choice predictor sizes shifted one
# column too far to the left (the BTB Entries value from column D was
# duplicated into column E, "Local Predictor Size"). The fix shifts the
# values one column to the right: D -> E -> F -> G, and leaves the now
# unused "BTB Entries" cell (D) blank, since BiMode predictors don't have
# a separate BTB-entries figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $predictor = $ws.Cells.Item($r, 3).Value2
    if ($predictor -eq "BiMode") {
        $dVal = $ws.Cells.Item($r, 4).Value2
        $eVal = $ws.Cells.Item($r, 5).Value2
        $fVal = $ws.Cells.Item($r, 6).Value2

        $ws.Cells.Item($r, 7).Value2 = $fVal
        $ws.Cells.Item($r, 6).Value2 = $eVal
        $ws.Cells.Item($r, 5).Value2 = $dVal
        $ws.Cells.Item($r, 4).ClearContents()
    }
}
